$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the (accidental) bold/font style that had been applied to C5,
# reverting the cell to the default (unstyled) format.
$ws.Range("C5").ClearFormats() | Out-Null

# Add a new "Background Sample Conversion" scenario row (row 13), following
# the same pattern as the existing rows 10-12.
$ws.Range("B13").Value = "Background Sample Conversion"
$ws.Range("C13").Value = "BackgroundSample.xlsx"
$ws.Range("D13").Formula = "=SUBSTITUTE(C13, "".xlsx"", "".ExpectedFeature"")"

# Match the selection left by the author after adding the new row.
$ws.Range("D12:D13").Select() | Out-Null
